# Append three new paragraphs at the end of the document (after the last,
# pre-existing empty paragraph, right before the sectPr). The first two are
# empty "terminal console" style paragraphs, the third carries a short note
# about mariadb installation problems.

$d = $word.ActiveDocument

# Build the shared <w:tabs> block used by this section's paragraphs
# (16 left tabs every 916 twips, matching the surrounding "terminal" text).
$tabPositions = 916,1832,2748,3664,4580,5496,6412,7328,8244,9160,10076,10992,11908,12824,13740,14656
$tabsXml = ""
foreach ($pos in $tabPositions) {
    $tabsXml += '<w:tab w:val="left" w:pos="' + $pos + '"/>'
}

# rPr shared by the paragraph mark of the two empty paragraphs, and by each
# run in the third (text) paragraph.
$rPrPlain = '<w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr>'

# rPr for the paragraph mark of the third paragraph (adds underline).
$rPrUnderline = '<w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:kern w:val="0"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="single"/><w:lang w:eastAsia="es-ES"/><w14:ligatures w14:val="none"/></w:rPr>'

$pPrCommon = '<w:pPr><w:tabs>' + $tabsXml + '</w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' + $rPrPlain + '</w:pPr>'
$pPrUnderline = '<w:pPr><w:tabs>' + $tabsXml + '</w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' + $rPrUnderline + '</w:pPr>'

# --- Paragraph 1: empty --------------------------------------------------
$para1Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $pPrCommon + '</w:p>'

# --- Paragraph 2: empty --------------------------------------------------
$para2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $pPrCommon + '</w:p>'

# --- Paragraph 3: the mariadb note, four runs -----------------------------
$run1 = '<w:r>' + $rPrPlain + '<w:t>Me estoy encontrando con muchos problemas para instalar mariadb.</w:t></w:r>'
$run2 = '<w:r>' + $rPrPlain + '<w:t xml:space="preserve"> He tenido que </w:t></w:r>'
$run3 = '<w:r>' + $rPrPlain + '<w:t xml:space="preserve">crear un usuario y darle los permisos necesarios en varias carpetas </w:t></w:r>'
$run4 = '<w:r>' + $rPrPlain + '<w:lastRenderedPageBreak/><w:t>para poder ejecutar el comando mysqld, ya que systemctl en Docker no funciona bien.</w:t></w:r>'
$para3Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + $pPrUnderline + $run1 + $run2 + $run3 + $run4 + '</w:p>'

# Insert the three paragraphs, one at a time, right after the document's
# current last paragraph, then replace each freshly-created paragraph's
# content with the exact XML we want (this avoids leaving a stray empty
# <w:r> that Word's COM InsertParagraphAfter would otherwise add).
$countBefore = $d.Paragraphs.Count

$last = $d.Paragraphs.Last
$null = $last.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$null = $p1.Range.InsertXML($para1Xml)

$p1 = $d.Paragraphs.Last
$null = $p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$null = $p2.Range.InsertXML($para2Xml)

$p2 = $d.Paragraphs.Last
$null = $p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$null = $p3.Range.InsertXML($para3Xml)

Write-Host "Paragraphs:" $countBefore "->" $d.Paragraphs.Count
